$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 3 data rows (old rows 2-4); remaining rows shift up
$ws.Range("A2:C4").EntireRow.Delete()

# Append 13 new data rows at the bottom (new rows 19-31)
$ws.Range("A19").Value = 2.29803731803496
$ws.Range("B19").Value = -5.991696847046798
$ws.Range("C19").Value = -5.515202137812292
$ws.Range("A20").Value = 2.744500100300473
$ws.Range("B20").Value = -0.9322929135791883
$ws.Range("C20").Value = -6.360266283544435
$ws.Range("A21").Value = -5.222654512415387
$ws.Range("B21").Value = -0.6178251758451783
$ws.Range("C21").Value = 4.528362251701193
$ws.Range("A22").Value = -5.398936099406894
$ws.Range("B22").Value = 0.6328901628238506
$ws.Range("C22").Value = 6.700743228821124
$ws.Range("A23").Value = -2.420953338994997
$ws.Range("B23").Value = 9.947797151136147
$ws.Range("C23").Value = -0.9126796104522215
$ws.Range("A24").Value = -0.05836680110209574
$ws.Range("B24").Value = 8.935851176995778
$ws.Range("C24").Value = -8.390108103527448
$ws.Range("A25").Value = 4.732191570142168
$ws.Range("B25").Value = -8.935019577985051
$ws.Range("C25").Value = -8.243645213661406
$ws.Range("A26").Value = 3.33527006372728
$ws.Range("B26").Value = -10.98798226935694
$ws.Range("C26").Value = -9.769438461483361
$ws.Range("A27").Value = -0.9925953642860486
$ws.Range("B27").Value = -2.305789720325507
$ws.Range("C27").Value = 3.166530378201873
$ws.Range("A28").Value = -4.912634885748004
$ws.Range("B28").Value = 4.700535394758443
$ws.Range("C28").Value = 3.522403849981169
$ws.Range("A29").Value = -3.346227371255714
$ws.Range("B29").Value = 3.616383455186579
$ws.Range("C29").Value = -1.352327122114249
$ws.Range("A30").Value = 3.400979183107231
$ws.Range("B30").Value = 13.99426472249469
$ws.Range("C30").Value = -1.426351852441958
$ws.Range("A31").Value = 0.4026054789882916
$ws.Range("B31").Value = -6.858872019183412
$ws.Range("C31").Value = -4.158537519539816

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
